$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Rushing": Week 17 rushing attempts/etc. logged for a few players
# ---------------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

$rushing.Range("E2").Value = 10   # T.Tagovailoa RZATT
$rushing.Range("C4").Value = 93   # M.Gaskin 1DATT
$rushing.Range("D4").Value = 64   # M.Gaskin 2DATT
$rushing.Range("C8").Value = 25   # D.Johnson 1DATT
$rushing.Range("D8").Value = 16   # D.Johnson 2DATT
$rushing.Range("E8").Value = 6    # D.Johnson RZATT
$rushing.Range("C9").Value = 15   # P.Lindsay 1DATT

# ---------------------------------------------------------------------------
# Sheet "Receiving": Week 17 data logged, including a brand-new row for
# P.Lindsay (inserted between D.Johnson and D.Parker).
# ---------------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# D.Johnson picked up Week 17 receiving stats too.
$receiving.Range("C6").Value = 3
$receiving.Range("D6").Value = 2

# Insert a new row 7, shifting D.Parker..D.Smythe down one row each.
$receiving.Rows("7:7").Insert()

# Pick up the same formatting (bold index column w/ border) as the row above.
$receiving.Range("A6:H6").Copy()
$receiving.Range("A7:H7").PasteSpecial(-4122)

# Fill in the new P.Lindsay row.
$receiving.Range("A7").Value = 5
$receiving.Range("B7").Value = "P.Lindsay"
$receiving.Range("C7").Value = 1
$receiving.Range("D7").Value = 1
$receiving.Range("E7").Value = 0
$receiving.Range("F7").Value = 0
$receiving.Range("G7").Value = 0
$receiving.Range("H7").Value = 0

# Re-sequence the index column for the rows that shifted down.
$receiving.Range("A8").Value = 6
$receiving.Range("A9").Value = 7
$receiving.Range("A10").Value = 8
$receiving.Range("A11").Value = 9
$receiving.Range("A12").Value = 10
$receiving.Range("A13").Value = 11
$receiving.Range("A14").Value = 12
$receiving.Range("A15").Value = 13
$receiving.Range("A16").Value = 14
$receiving.Range("A17").Value = 15
$receiving.Range("A18").Value = 16
$receiving.Range("A19").Value = 17
$receiving.Range("A20").Value = 18

# Week 17 stat updates for the players who played that week.
$receiving.Range("C8").Value = 86   # D.Parker Short Target
$receiving.Range("D8").Value = 70   # D.Parker Short Comp
$receiving.Range("E8").Value = 19   # D.Parker Deep Target
$receiving.Range("F8").Value = 9    # D.Parker Deep Comp

$receiving.Range("C9").Value = 123  # J.Waddle Short Target
$receiving.Range("D9").Value = 99   # J.Waddle Short Comp
$receiving.Range("E9").Value = 23   # J.Waddle Deep Target
$receiving.Range("F9").Value = 11   # J.Waddle Deep Comp

$receiving.Range("C13").Value = 18  # M.Hollins Short Target

$receiving.Range("C14").Value = 12  # I.Ford Short Target
$receiving.Range("D14").Value = 9   # I.Ford Short Comp

$receiving.Range("C17").Value = 87  # M.Gesicki Short Target
$receiving.Range("D17").Value = 59  # M.Gesicki Short Comp
$receiving.Range("E17").Value = 23  # M.Gesicki Deep Target

$receiving.Range("C20").Value = 32  # D.Smythe Short Target
$receiving.Range("D20").Value = 27  # D.Smythe Short Comp
$receiving.Range("E20").Value = 7   # D.Smythe Deep Target
$receiving.Range("F20").Value = 5   # D.Smythe Deep Comp
